$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain stored as text (matching original inlineStr cells)
$textCells = @('D5', 'D6', 'D7', 'D8', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D18', 'D19', 'D20', 'D21', 'D23', 'D27', 'D28', 'D29', 'D30', 'D31', 'D34', 'D35', 'D36', 'D37', 'D38', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D49', 'D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values from the crypto price refresh
$ws.Range('D2').Value = '24.544.57'
$ws.Range('E2').Value = '  +3.27%  '
$ws.Range('D3').Value = '1.693.03'
$ws.Range('E3').Value = '  +1.70%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').Value = '315.71'
$ws.Range('E5').Value = '  +2.08%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('D7').Value = '0.3932'
$ws.Range('E7').Value = '  +1.38%  '
$ws.Range('D8').Value = '0.4002'
$ws.Range('E8').Value = '  +1.48%  '
$ws.Range('E9').Value = '  +4.89%  '
$ws.Range('D10').Value = '1.002'
$ws.Range('E10').Value = '  +0.30%  '
$ws.Range('D11').Value = '53.15'
$ws.Range('E11').Value = '  +5.54%  '
$ws.Range('D12').Value = '0.08724'
$ws.Range('E12').Value = '  +0.80%  '
$ws.Range('D13').Value = '7.182'
$ws.Range('E13').Value = '  +6.84%  '
$ws.Range('D14').Value = '23.08'
$ws.Range('E14').Value = '  +2.10%  '
$ws.Range('D15').Value = '0.00001314'
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('D16').Value = '7.584'
$ws.Range('E16').Value = '  +4.13%  '
$ws.Range('D17').Value = '1.701.24'
$ws.Range('E17').Value = '  +1.81%  '
$ws.Range('D18').Value = '99.67'
$ws.Range('E18').Value = '  -0.03%  '
$ws.Range('D19').Value = '0.07044'
$ws.Range('E19').Value = '  +3.73%  '
$ws.Range('D20').Value = '19.58'
$ws.Range('E20').Value = '  +2.65%  '
$ws.Range('D21').Value = '6.831'
$ws.Range('E21').Value = '  +3.10%  '
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').Value = '14.02'
$ws.Range('E23').Value = '  +1.28%  '
$ws.Range('D24').Value = '24.533.48'
$ws.Range('E24').Value = '  +3.28%  '
$ws.Range('E25').Value = '  +6.56%  '
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('D27').Value = '22.31'
$ws.Range('E27').Value = '  +2.81%  '
$ws.Range('D28').Value = '160.69'
$ws.Range('E28').Value = '  +0.84%  '
$ws.Range('D29').Value = '5.211'
$ws.Range('E29').Value = '  +0.25%  '
$ws.Range('D30').Value = '134.21'
$ws.Range('E30').Value = '  +3.76%  '
$ws.Range('D31').Value = '7.465'
$ws.Range('E31').Value = '  +12.08%  '
$ws.Range('D32').Value = '1.882.79'
$ws.Range('E32').Value = '  +1.52%  '
$ws.Range('E33').Value = '  -2.83%  '
$ws.Range('D34').Value = '0.08521'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').Value = '7.264'
$ws.Range('E35').Value = '  +10.18%  '
$ws.Range('D36').Value = '11.35'
$ws.Range('E36').Value = '  +8.44%  '
$ws.Range('D37').Value = '1.957'
$ws.Range('E37').Value = '  -1.47%  '
$ws.Range('D38').Value = '0.2706'
$ws.Range('E38').Value = '  +1.85%  '
$ws.Range('E39').Value = '  -0.33%  '
$ws.Range('D40').Value = '0.02745'
$ws.Range('E40').Value = '  +9.19%  '
$ws.Range('D41').Value = '0.09022'
$ws.Range('E41').Value = '  +2.78%  '
$ws.Range('D42').Value = '1.474'
$ws.Range('E42').Value = '  +1.24%  '
$ws.Range('D43').Value = '0.7621'
$ws.Range('E43').Value = '  +1.10%  '
$ws.Range('D44').Value = '0.7172'
$ws.Range('E44').Value = '  +2.24%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').Value = '2.524'
$ws.Range('E45').Value = '  +4.74%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '15.32'
$ws.Range('E46').Value = '  +3.56%  '
$ws.Range('D47').Value = '4.202'
$ws.Range('E47').Value = '  +2.49%  '
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('D49').Value = '140.85'
$ws.Range('E49').Value = '  +1.49%  '
$ws.Range('D50').Value = '1.307'
$ws.Range('E50').Value = '  +4.49%  '
